{"js": "// UPOTREBLJENA METODOLOGIJA U OVOM PROJEKTU\n//\n// 1) Remove the stray empty centered paragraph that sits right after the\n//    very first (header) table in the document body.\n// 2) Add a new row (RB \"65\") to the \"Primjedbe na projekt i dizajn\" table,\n//    right after the existing last row (\"64\" / TOMISLAV ZDUNI\u0106).\n// 3) Because content shifted earlier in the document, Word's cached\n//    \"last rendered page break\" marker moves from the run holding \"32\"\n//    to the run holding \"34\" (two rows later) in the same table.\n\nconst body = context.document.body;\n\n// --- 1. Delete the empty paragraph right after the first table ---------\n// `body.paragraphs` is a flattened list that also walks into table cells,\n// so look for the first *body-level* (not inside any table) empty\n// paragraph - that is the stray paragraph sitting right after the first\n// (header) table, before the \"ZAPISNIK SA KONZULTACIJE...\" title.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst parentTables = paragraphs.items.map((p) => p.parentTableOrNullObject);\nparentTables.forEach((t) => t.load(\"isNullObject\"));\nawait context.sync();\n\nlet targetPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\" && parentTables[i].isNullObject) {\n    targetPara = paragraphs.items[i];\n    break;\n  }\n}\nif (targetPara) {\n  targetPara.delete();\n  await context.sync();\n}\n\n// --- 2. Append the new \"65\" row to the remarks table --------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The remarks table (\"Primjedbe na projekt i dizajn\") is the second table\n// in the document and ends with row RB \"64\" / \"TOMISLAV ZDUNI\u0106\".\nconst remarksTable = tables.items[1];\nremarksTable.addRows(\"End\", 1, [\n  [\"65\", \"Tekst, dijagrami, ostalo\", \"Obavezno navesti izvor\", \"SVI\"]\n]);\nawait context.sync();\n\n// --- 3. Move <w:lastRenderedPageBreak/> from run \"32\" to run \"34\" ------\nconst ooxmlNoBreak = (num) =>\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:t>' + num + '</w:t></w:r></w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nconst ooxmlWithBreak = (num) =>\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>' + num + '</w:t></w:r></w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nconst results32 = body.search(\"32\", { matchWholeWord: true });\nresults32.load(\"items\");\nawait context.sync();\nif (results32.items.length > 0) {\n  results32.items[0].insertOoxml(ooxmlNoBreak(\"32\"), \"Replace\");\n  await context.sync();\n}\n\nconst results34 = body.search(\"34\", { matchWholeWord: true });\nresults34.load(\"items\");\nawait context.sync();\nif (results34.items.length > 0) {\n  results34.items[0].insertOoxml(ooxmlWithBreak(\"34\"), \"Replace\");\n  await context.sync();\n}\n", "ps1": "# UPOTREBLJENA METODOLOGIJA U OVOM PROJEKTU\n#\n# 1) Remove the stray empty centered paragraph that sits right after the\n#    very first (header) table in the document body.\n# 2) Add a new row (RB \"65\") to the \"Primjedbe na projekt i dizajn\" table,\n#    right after the existing last row (\"64\" / TOMISLAV ZDUNIC).\n# 3) Because content shifted earlier in the document, Word's cached\n#    \"last rendered page break\" marker moves from the run holding \"32\"\n#    to the run holding \"34\" (two rows later) in the same table.\n\n$d = $word.ActiveDocument\n\n# --- 1. Delete the empty paragraph right after the first table ---------\n$firstTable = $d.Tables.Item(1)\n$posAfterTable = $firstTable.Range.End\n$strayParaRange = $d.Range($posAfterTable, $posAfterTable + 1)\n$strayParaRange.Delete()\n\n# --- 2. Append the new \"65\" row to the remarks table --------------------\n# The remarks table (\"Primjedbe na projekt i dizajn\") is the second table\n# in the document and ends with row RB \"64\" / \"TOMISLAV ZDUNIC\".\n$remarksTable = $d.Tables.Item(2)\n$newRow = $remarksTable.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"65\"\n$newRow.Cells.Item(2).Range.Text = \"Tekst, dijagrami, ostalo\"\n$newRow.Cells.Item(3).Range.Text = \"Obavezno navesti izvor\"\n$newRow.Cells.Item(4).Range.Text = \"SVI\"\n\n# --- 3. Move <w:lastRenderedPageBreak/> from run \"32\" to run \"34\" ------\n$pPrXml = '<w:pPr><w:ind w:firstLine=\"0\"/><w:jc w:val=\"left\"/></w:pPr>'\n\n$rng32 = $d.Content\n$rng32.Find.ClearFormatting()\n$rng32.Find.Text = \"32\"\n$rng32.Find.MatchWholeWord = $true\n$rng32.Find.Execute() | Out-Null\nif ($rng32.Find.Found) {\n    $ooxmlNoBreak = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $pPrXml + '<w:r><w:t>32</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $rng32.InsertXML($ooxmlNoBreak)\n}\n\n$rng34 = $d.Content\n$rng34.Find.ClearFormatting()\n$rng34.Find.Text = \"34\"\n$rng34.Find.MatchWholeWord = $true\n$rng34.Find.Execute() | Out-Null\nif ($rng34.Find.Found) {\n    $ooxmlWithBreak = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $pPrXml + '<w:r><w:lastRenderedPageBreak/><w:t>34</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $rng34.InsertXML($ooxmlWithBreak)\n}\n"}
